$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting so
# values such as "1.001" or "0.06100" are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.365.03"
$ws.Range("E2").Value = "  +3.23%  "

# Row 3
$ws.Range("D3").Value = "1.868.83"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "339.06"
$ws.Range("E5").Value = "  +1.90%  "

# Row 6
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").Value = "0.4704"
$ws.Range("E7").Value = "  +2.01%  "

# Row 8
$ws.Range("E8").Value = "  +3.67%  "

# Row 9
$ws.Range("D9").Value = "47.54"
$ws.Range("E9").Value = "  +2.46%  "

# Row 10
$ws.Range("D10").Value = "0.08042"
$ws.Range("E10").Value = "  +1.69%  "

# Row 11
$ws.Range("D11").Value = "1.001"
$ws.Range("E11").Value = "  +2.50%  "

# Row 12
$ws.Range("D12").Value = "21.93"
$ws.Range("E12").Value = "  +3.52%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.881.52"
$ws.Range("E13").Value = "  +1.00%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.041"
$ws.Range("E14").Value = "  +2.24%  "

# Row 15
$ws.Range("D15").Value = "7.253"
$ws.Range("E15").Value = "  +2.79%  "

# Row 16
$ws.Range("D16").Value = "91.22"
$ws.Range("E16").Value = "  +3.59%  "

# Row 17
$ws.Range("E17").Value = "  -0.12%  "

# Row 18
$ws.Range("D18").Value = "0.00001044"
$ws.Range("E18").Value = "  +1.35%  "

# Row 19
$ws.Range("D19").Value = "0.06622"
$ws.Range("E19").Value = "  -0.27%  "

# Row 20
$ws.Range("D20").Value = "17.56"
$ws.Range("E20").Value = "  +2.85%  "

# Row 21
$ws.Range("E21").Value = "  -0.24%  "

# Row 22
$ws.Range("D22").Value = "28.379.21"
$ws.Range("E22").Value = "  +3.31%  "

# Row 23
$ws.Range("D23").Value = "5.479"
$ws.Range("E23").Value = "  +1.99%  "

# Row 24
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +1.74%  "

# Row 25
$ws.Range("E25").Value = "  -1.99%  "

# Row 26
$ws.Range("D26").Value = "2.093.95"
$ws.Range("E26").Value = "  +0.67%  "

# Row 27
$ws.Range("D27").Value = "161.01"
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
$ws.Range("E28").Value = "  +1.73%  "

# Row 29
$ws.Range("D29").Value = "2.123"
$ws.Range("E29").Value = "  +2.52%  "

# Row 30
$ws.Range("D30").Value = "5.508"
$ws.Range("E30").Value = "  +3.20%  "

# Row 31
$ws.Range("D31").Value = "120.39"
$ws.Range("E31").Value = "  +1.12%  "

# Row 32
$ws.Range("D32").Value = "0.9707"
$ws.Range("E32").Value = "  +1.43%  "

# Row 33
$ws.Range("D33").Value = "0.09511"
$ws.Range("E33").Value = "  +2.25%  "

# Row 34
$ws.Range("D34").Value = "3.596"
$ws.Range("E34").Value = "  +0.73%  "

# Row 35
$ws.Range("D35").Value = "5.350"
$ws.Range("E35").Value = "  +1.87%  "

# Row 36
$ws.Range("D36").Value = "1.374"
$ws.Range("E36").Value = "  +3.90%  "

# Row 37
$ws.Range("D37").Value = "0.06100"
$ws.Range("E37").Value = "  +2.52%  "

# Row 38
$ws.Range("D38").Value = "0.02255"
$ws.Range("E38").Value = "  +2.90%  "

# Row 39
$ws.Range("D39").Value = "8.385"
$ws.Range("E39").Value = "  +3.69%  "

# Row 40
$ws.Range("D40").Value = "1.188"

# Row 41
$ws.Range("D41").Value = "0.5957"
$ws.Range("E41").Value = "  +2.47%  "

# Row 42
$ws.Range("E42").Value = "  -0.19%  "

# Row 43
$ws.Range("D43").Value = "0.1875"
$ws.Range("E43").Value = "  +1.70%  "

# Row 44
$ws.Range("D44").Value = "10.38"
$ws.Range("E44").Value = "  +3.28%  "

# Row 45
$ws.Range("D45").Value = "1.288"
$ws.Range("E45").Value = "  +3.08%  "

# Row 46
$ws.Range("D46").Value = "0.5595"
$ws.Range("E46").Value = "  +1.73%  "

# Row 47
$ws.Range("D47").Value = "12.22"
$ws.Range("E47").Value = "  +2.04%  "

# Row 48
$ws.Range("D48").Value = "1.958"
$ws.Range("E48").Value = "  +4.52%  "

# Row 49
$ws.Range("D49").Value = "0.06882"
$ws.Range("E49").Value = "  +3.39%  "

# Row 50
$ws.Range("E50").Value = "  +17.44%  "

# Row 51
$ws.Range("D51").Value = "111.63"
$ws.Range("E51").Value = "  +1.37%  "
